$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wb.Unprotect() | Out-Null
$ws.Range("A1").Value = "2023-07-23 00:00:00"
$ws.Range("B1").Value = 2
$ws.Range("A2:B2").Select() | Out-Null
